$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = '28.281.02'
$ws.Cells.Item(2,5).Value = '  -5.53%  '

$ws.Cells.Item(3,4).Value = '1.838.82'
$ws.Cells.Item(3,5).Value = '  -5.17%  '

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = '1.002'
$ws.Cells.Item(4,4).Style = "Normal"
$ws.Cells.Item(4,5).Value = '  -0.42%  '

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = '330.45'
$ws.Cells.Item(5,4).Style = "Normal"
$ws.Cells.Item(5,5).Value = '  -1.59%  '

$ws.Cells.Item(6,5).Value = '  -0.27%  '

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = '0.4602'
$ws.Cells.Item(7,4).Style = "Normal"
$ws.Cells.Item(7,5).Value = '  -4.90%  '

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = '0.3859'
$ws.Cells.Item(8,4).Style = "Normal"
$ws.Cells.Item(8,5).Value = '  -6.39%  '

$ws.Cells.Item(9,2).Value = 'Dogecoin'
$ws.Cells.Item(9,3).Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '0.07869'
$ws.Cells.Item(9,4).Style = "Normal"
$ws.Cells.Item(9,5).Value = '  -3.71%  '

$ws.Cells.Item(10,2).Value = 'Polygon'
$ws.Cells.Item(10,3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0.9659'
$ws.Cells.Item(10,4).Style = "Normal"
$ws.Cells.Item(10,5).Value = '  -4.99%  '

$ws.Cells.Item(11,2).Value = 'Solana'
$ws.Cells.Item(11,3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = '21.96'
$ws.Cells.Item(11,4).Style = "Normal"
$ws.Cells.Item(11,5).Value = '  -7.25%  '

$ws.Cells.Item(12,2).Value = 'WrappedEther'
$ws.Cells.Item(12,3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(12,4).Value = '1.861.01'
$ws.Cells.Item(12,5).Value = '  -4.49%  '

$ws.Cells.Item(13,2).Value = 'Polkadot'
$ws.Cells.Item(13,3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '5.718'
$ws.Cells.Item(13,4).Style = "Normal"
$ws.Cells.Item(13,5).Value = '  -6.13%  '

$ws.Cells.Item(14,2).Value = 'Chainlink'
$ws.Cells.Item(14,3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '6.921'
$ws.Cells.Item(14,4).Style = "Normal"
$ws.Cells.Item(14,5).Value = '  -5.16%  '

$ws.Cells.Item(15,2).Value = 'TRON'
$ws.Cells.Item(15,3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0.06888'
$ws.Cells.Item(15,4).Style = "Normal"
$ws.Cells.Item(15,5).Value = '  +0.50%  '

$ws.Cells.Item(16,2).Value = 'BinanceUSD'
$ws.Cells.Item(16,3).Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = '1.002'
$ws.Cells.Item(16,4).Style = "Normal"
$ws.Cells.Item(16,5).Value = '  -0.39%  '

$ws.Cells.Item(17,2).Value = 'Litecoin'
$ws.Cells.Item(17,3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = '86.89'
$ws.Cells.Item(17,4).Style = "Normal"
$ws.Cells.Item(17,5).Value = '  -4.71%  '

$ws.Cells.Item(18,2).Value = 'ShibaInu'
$ws.Cells.Item(18,3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = '0.000009954'
$ws.Cells.Item(18,4).Style = "Normal"
$ws.Cells.Item(18,5).Value = '  -3.95%  '

$ws.Cells.Item(19,2).Value = 'Avalanche'
$ws.Cells.Item(19,3).Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = '16.93'
$ws.Cells.Item(19,4).Style = "Normal"
$ws.Cells.Item(19,5).Value = '  -4.90%  '

$ws.Cells.Item(20,2).Value = 'Dai'
$ws.Cells.Item(20,3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value = '1.003'
$ws.Cells.Item(20,4).Style = "Normal"
$ws.Cells.Item(20,5).Value = '  -0.23%  '

$ws.Cells.Item(21,2).Value = 'WrappedBTC'
$ws.Cells.Item(21,3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(21,4).Value = '28.302.93'
$ws.Cells.Item(21,5).Value = '  -5.42%  '

$ws.Cells.Item(22,2).Value = 'Uniswap'
$ws.Cells.Item(22,3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = '5.342'
$ws.Cells.Item(22,4).Style = "Normal"
$ws.Cells.Item(22,5).Value = '  -5.20%  '

$ws.Cells.Item(23,2).Value = 'Cosmos'
$ws.Cells.Item(23,3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = '10.98'
$ws.Cells.Item(23,4).Style = "Normal"
$ws.Cells.Item(23,5).Value = '  -7.49%  '

$ws.Cells.Item(24,2).Value = 'Toncoin'
$ws.Cells.Item(24,3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = '2.149'
$ws.Cells.Item(24,4).Style = "Normal"
$ws.Cells.Item(24,5).Value = '  -1.57%  '

$ws.Cells.Item(25,2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(25,3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(25,4).Value = '2.081.38'
$ws.Cells.Item(25,5).Value = '  -5.02%  '

$ws.Cells.Item(26,2).Value = 'Monero'
$ws.Cells.Item(26,3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(26,4).NumberFormat = "@"
$ws.Cells.Item(26,4).Value = '153.73'
$ws.Cells.Item(26,4).Style = "Normal"
$ws.Cells.Item(26,5).Value = '  -1.81%  '

$ws.Cells.Item(27,2).Value = 'EthereumClassic'
$ws.Cells.Item(27,3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = '19.20'
$ws.Cells.Item(27,4).Style = "Normal"
$ws.Cells.Item(27,5).Value = '  -4.43%  '

$ws.Cells.Item(28,2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(28,3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = '5.767'
$ws.Cells.Item(28,4).Style = "Normal"
$ws.Cells.Item(28,5).Value = '  -13.84%  '

$ws.Cells.Item(29,2).Value = 'LidoDAOToken'
$ws.Cells.Item(29,3).Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = '1.985'
$ws.Cells.Item(29,4).Style = "Normal"
$ws.Cells.Item(29,5).Value = '  -5.47%  '

$ws.Cells.Item(30,2).Value = 'BitcoinCash'
$ws.Cells.Item(30,3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = '117.06'
$ws.Cells.Item(30,4).Style = "Normal"
$ws.Cells.Item(30,5).Value = '  -3.68%  '

$ws.Cells.Item(31,2).Value = 'ImmutableX'
$ws.Cells.Item(31,3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = '0.9419'
$ws.Cells.Item(31,4).Style = "Normal"
$ws.Cells.Item(31,5).Value = '  -6.71%  '

$ws.Cells.Item(32,2).Value = 'Stellar'
$ws.Cells.Item(32,3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = '0.09311'
$ws.Cells.Item(32,4).Style = "Normal"
$ws.Cells.Item(32,5).Value = '  -3.34%  '

$ws.Cells.Item(33,2).Value = 'Filecoin'
$ws.Cells.Item(33,3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = '5.288'
$ws.Cells.Item(33,4).Style = "Normal"
$ws.Cells.Item(33,5).Value = '  -5.37%  '

$ws.Cells.Item(34,2).Value = 'HuobiToken'
$ws.Cells.Item(34,3).Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = '3.447'
$ws.Cells.Item(34,4).Style = "Normal"
$ws.Cells.Item(34,5).Value = '  -2.74%  '

$ws.Cells.Item(35,2).Value = 'ARBITRUM'
$ws.Cells.Item(35,3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(35,4).NumberFormat = "@"
$ws.Cells.Item(35,4).Value = '1.327'
$ws.Cells.Item(35,4).Style = "Normal"
$ws.Cells.Item(35,5).Value = '  -6.92%  '

$ws.Cells.Item(36,2).Value = 'Hedera'
$ws.Cells.Item(36,3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = '0.06015'
$ws.Cells.Item(36,4).Style = "Normal"
$ws.Cells.Item(36,5).Value = '  -8.25%  '

$ws.Cells.Item(37,2).Value = 'VeChain'
$ws.Cells.Item(37,3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(37,4).NumberFormat = "@"
$ws.Cells.Item(37,4).Value = '0.02153'
$ws.Cells.Item(37,4).Style = "Normal"
$ws.Cells.Item(37,5).Value = '  -5.93%  '

$ws.Cells.Item(38,2).Value = 'TrustWalletToken'
$ws.Cells.Item(38,3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = '1.150'
$ws.Cells.Item(38,4).Style = "Normal"
$ws.Cells.Item(38,5).Value = '  -4.84%  '

$ws.Cells.Item(39,2).Value = 'Frax'
$ws.Cells.Item(39,3).Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = '1.001'
$ws.Cells.Item(39,4).Style = "Normal"
$ws.Cells.Item(39,5).Value = '  -0.28%  '

$ws.Cells.Item(40,2).Value = 'FraxShare'
$ws.Cells.Item(40,3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = '7.614'
$ws.Cells.Item(40,4).Style = "Normal"
$ws.Cells.Item(40,5).Value = '  -4.73%  '

$ws.Cells.Item(41,2).Value = 'TheSandbox'
$ws.Cells.Item(41,3).Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value = '0.5629'
$ws.Cells.Item(41,4).Style = "Normal"
$ws.Cells.Item(41,5).Value = '  -5.78%  '

$ws.Cells.Item(42,2).Value = 'Aptos'
$ws.Cells.Item(42,3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = '10.02'
$ws.Cells.Item(42,4).Style = "Normal"
$ws.Cells.Item(42,5).Value = '  -6.72%  '

$ws.Cells.Item(43,2).Value = 'Algorand'
$ws.Cells.Item(43,3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = '0.1784'
$ws.Cells.Item(43,4).Style = "Normal"
$ws.Cells.Item(43,5).Value = '  -3.51%  '

$ws.Cells.Item(44,2).Value = 'WEMIXToken'
$ws.Cells.Item(44,3).Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = '1.237'
$ws.Cells.Item(44,4).Style = "Normal"
$ws.Cells.Item(44,5).Value = '  -2.92%  '

$ws.Cells.Item(45,2).Value = 'RenderToken'
$ws.Cells.Item(45,3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = '2.278'
$ws.Cells.Item(45,4).Style = "Normal"
$ws.Cells.Item(45,5).Value = '  -9.11%  '

$ws.Cells.Item(46,2).Value = 'EnergySwap'
$ws.Cells.Item(46,3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value = '11.67'
$ws.Cells.Item(46,4).Style = "Normal"
$ws.Cells.Item(46,5).Value = '  -5.52%  '

$ws.Cells.Item(47,2).Value = 'Decentraland'
$ws.Cells.Item(47,3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = '0.5283'
$ws.Cells.Item(47,4).Style = "Normal"
$ws.Cells.Item(47,5).Value = '  -5.27%  '

$ws.Cells.Item(48,2).Value = 'Cronos'
$ws.Cells.Item(48,3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(48,4).NumberFormat = "@"
$ws.Cells.Item(48,4).Value = '0.07027'
$ws.Cells.Item(48,4).Style = "Normal"
$ws.Cells.Item(48,5).Value = '  -6.20%  '

$ws.Cells.Item(49,2).Value = 'NEARProtocol'
$ws.Cells.Item(49,3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(49,4).NumberFormat = "@"
$ws.Cells.Item(49,4).Value = '1.839'
$ws.Cells.Item(49,4).Style = "Normal"
$ws.Cells.Item(49,5).Value = '  -7.70%  '

$ws.Cells.Item(50,2).Value = 'Quant'
$ws.Cells.Item(50,3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = '112.86'
$ws.Cells.Item(50,4).Style = "Normal"
$ws.Cells.Item(50,5).Value = '  -3.74%  '

$ws.Cells.Item(51,2).Value = 'PaxDollar'
$ws.Cells.Item(51,3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = '1.002'
$ws.Cells.Item(51,4).Style = "Normal"
$ws.Cells.Item(51,5).Value = '  -0.54%  '
